$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add two new rows (18,19) comparing Tlin/Hsou for Ra=10^8
$ws.Range("C18").Value = "Ra=10^8"
$ws.Range("D18").Value = "Tlin"
$ws.Range("C19").Value = "Ra=10^8"
$ws.Range("D19").Value = "Hsou"

# Update the grid-dependence test labels (rows 7-14): 80x80/120x120/160x160/200x200 -> 100x100/140x140/180x180/225x225
$ws.Range("C8").Value  = "140x140"
$ws.Range("C9").Value  = "180x180"
$ws.Range("C10").Value = "225x225"
$ws.Range("C7").Value  = "100x100"
$ws.Range("C11").Value = "100x100"
$ws.Range("C12").Value = "140x140"
$ws.Range("C13").Value = "180x180"
$ws.Range("C14").Value = "225x225"

# Update selection to match the edited workbook state
$ws.Range("F13").Select()
